$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28. This shifts the existing row 28
# (M. Radaszkiewicz) down to row 29, and row 29 (K. Moskal) down to row 30.
$ws.Rows.Item(28).Insert()

# Copy the formatting of the row above (row 27) into the new row 28 so that
# column A keeps its bold / bordered / centered style.
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)

# Fill the newly inserted row 28 with S. Jurić's stats.
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "S. Jurić"
$ws.Range("C28").Value = "'23"
$ws.Range("D28").Value = "N"
$ws.Range("E28").Value = "'0"
$ws.Range("F28").Value = "'0"
$ws.Range("G28").Value = "'0"
$ws.Range("H28").Value = "'0"
$ws.Range("I28").Value = "'0"
$ws.Range("J28").Value = "'0"
$ws.Range("K28").Value = "'0"
$ws.Range("L28").Value = "'0"
$ws.Range("M28").Value = "'0"
$ws.Range("N28").Value = "'0"

# The row-index column (A) is a manually maintained sequence (row - 2), so
# the rows that were pushed down need their index renumbered to stay
# sequential.
$ws.Range("A29").Value = 27
$ws.Range("A30").Value = 28
